# Update "想去人数" (number of people interested) values in the
# "展览" and "全部类型" worksheets to reflect the latest generated data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 82
$ws1.Range("F4").Value = 2213
$ws1.Range("F5").Value = 194
$ws1.Range("F6").Value = 367

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 82
$ws4.Range("F4").Value = 2213
$ws4.Range("F5").Value = 194
$ws4.Range("F7").Value = 367
